# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interest count) figures in column F across the
# three data sheets that carry this metric: 展览 (Exhibitions), 本地生活
# (Local life), and 全部类型 (All types, the combined sheet).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 750
$ws.Range("F6").Value = 2379
$ws.Range("F7").Value = 53
$ws.Range("F8").Value = 1765
$ws.Range("F9").Value = 2990
$ws.Range("F11").Value = 4455
$ws.Range("F12").Value = 392
$ws.Range("F15").Value = 560
$ws.Range("F17").Value = 617
$ws.Range("F18").Value = 223
$ws.Range("F21").Value = 307
$ws.Range("F22").Value = 4514
$ws.Range("F24").Value = 3934
$ws.Range("F25").Value = 1142
$ws.Range("F26").Value = 218
$ws.Range("F27").Value = 574
$ws.Range("F30").Value = 608
$ws.Range("F31").Value = 592
$ws.Range("F32").Value = 553

# --- Sheet 3: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 210
$ws.Range("F3").Value = 1039
$ws.Range("F4").Value = 19

# --- Sheet 4: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 210
$ws.Range("F4").Value = 1039
$ws.Range("F5").Value = 19
$ws.Range("F8").Value = 750
$ws.Range("F9").Value = 2379
$ws.Range("F10").Value = 53
$ws.Range("F11").Value = 1765
$ws.Range("F13").Value = 2990
$ws.Range("F15").Value = 4455
$ws.Range("F16").Value = 392
$ws.Range("F19").Value = 560
$ws.Range("F21").Value = 617
$ws.Range("F22").Value = 223
$ws.Range("F26").Value = 307
$ws.Range("F27").Value = 4514
$ws.Range("F29").Value = 3934
$ws.Range("F30").Value = 1142
$ws.Range("F31").Value = 218
$ws.Range("F32").Value = 574
$ws.Range("F35").Value = 608
$ws.Range("F36").Value = 592
$ws.Range("F37").Value = 553
